$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix capitalization of the class-name labels in column A (dropdown / non-dropdown
# classes). Only the text casing changes; the values in column B stay the same.
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

# Move the active selection to A8, matching the saved workbook state.
$ws.Range("A8").Select()
